$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 26322304
$ws.Range("I62").Value = 1418.6666
$ws.Range("J62").Value = 125025624
$ws.Range("K62").Value = 1418.6666
$ws.Range("L62").Value = 125025624
$ws.Range("M62").Value = -794.6666
$ws.Range("N62").Value = -125026872
$ws.Range("H65").Value = 26322304
$ws.Range("I65").Value = 1418.6666
$ws.Range("J65").Value = 125025624
$ws.Range("K65").Value = 7093.333000000001
$ws.Range("L65").Value = 625128120
$ws.Range("M65").Value = -3973.333000000001
$ws.Range("N65").Value = -625134360
$ws.Range("H88").Value = 2128.7
$ws.Range("J88").Value = 2276.3333
$ws.Range("L88").Value = 2276.3333
$ws.Range("N88").Value = -3088.3333
$ws.Range("H91").Value = 2128.7
$ws.Range("J91").Value = 2276.3333
$ws.Range("L91").Value = 2276.3333
$ws.Range("N91").Value = -5084.3333
$ws.Range("H100").Value = 19609258
$ws.Range("I100").Value = 2150.5
$ws.Range("J100").Value = 37037796
$ws.Range("K100").Value = 2150.5
$ws.Range("L100").Value = 37037796
$ws.Range("M100").Value = -1609.5
$ws.Range("N100").Value = -37038878
$ws.Range("H124").Value = 34540
$ws.Range("J124").Value = 34540
$ws.Range("L124").Value = 34540
$ws.Range("N124").Value = -44360
$ws.Range("H137").Value = 4168011.8
$ws.Range("I137").Value = 1924258.5
$ws.Range("J137").Value = 10001770
$ws.Range("K137").Value = 5772775.5
$ws.Range("L137").Value = 30005310
$ws.Range("M137").Value = -5770225.5
$ws.Range("N137").Value = -30010410
$ws.Range("H138").Value = 2353.5356
$ws.Range("I138").Value = 2421.647
$ws.Range("J138").Value = 2323.8462
$ws.Range("K138").Value = 7264.941
$ws.Range("L138").Value = 6971.5386
$ws.Range("M138").Value = -2124.941
$ws.Range("N138").Value = -17251.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21637.152
$ws.Range("I32").Value = 4298.983
$ws.Range("K32").Value = 4298.983
$ws.Range("M32").Value = -4011.983
$ws.Range("H33").Value = 500013020
$ws.Range("I33").Value = 26026
$ws.Range("K33").Value = 26026
$ws.Range("M33").Value = -25697
$ws.Range("H122").Value = 1604.6316
$ws.Range("I122").Value = 1378
$ws.Range("J122").Value = 2454.5
$ws.Range("K122").Value = 4134
$ws.Range("L122").Value = 7363.5
$ws.Range("M122").Value = -1684
$ws.Range("N122").Value = -12263.5
$ws.Range("H132").Value = 241822.77
$ws.Range("I132").Value = 418813.84
$ws.Range("J132").Value = 5834.6665
$ws.Range("K132").Value = 1256441.52
$ws.Range("L132").Value = 17503.9995
$ws.Range("M132").Value = -1253911.52
$ws.Range("N132").Value = -22563.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2449.6924
$ws.Range("J86").Value = 1990
$ws.Range("L86").Value = 1990
$ws.Range("N86").Value = -4236
$ws.Range("H89").Value = 2449.6924
$ws.Range("J89").Value = 1990
$ws.Range("L89").Value = 9950
$ws.Range("N89").Value = -21182
$ws.Range("H94").Value = 931.4
$ws.Range("I94").Value = 1048.3914
$ws.Range("K94").Value = 1048.3914
$ws.Range("M94").Value = -597.3914
$ws.Range("H99").Value = 1813.931
$ws.Range("I99").Value = 1159.5333
$ws.Range("J99").Value = 2515.0715
$ws.Range("K99").Value = 1159.5333
$ws.Range("L99").Value = 2515.0715
$ws.Range("M99").Value = 338.4666999999999
$ws.Range("N99").Value = -5511.0715
$ws.Range("H105").Value = 3303.75
$ws.Range("I105").Value = 2147.4119
$ws.Range("J105").Value = 5090.8184
$ws.Range("K105").Value = 2147.4119
$ws.Range("L105").Value = 5090.8184
$ws.Range("M105").Value = -400.4119000000001
$ws.Range("N105").Value = -8584.8184
$ws.Range("H124").Value = 46250
$ws.Range("J124").Value = 46250
$ws.Range("L124").Value = 46250
$ws.Range("N124").Value = -56070
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H134").Value = 54052.895
$ws.Range("I134").Value = 58024.434
$ws.Range("J134").Value = 1430
$ws.Range("K134").Value = 174073.302
$ws.Range("L134").Value = 4290
$ws.Range("M134").Value = -171538.302
$ws.Range("N134").Value = -9360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 613.55316
$ws.Range("I113").Value = 586.2143
$ws.Range("J113").Value = 653.8421
$ws.Range("K113").Value = 1758.6429
$ws.Range("L113").Value = 1961.5263
$ws.Range("M113").Value = 411.3571000000002
$ws.Range("N113").Value = -6301.5263

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 17765.666
$ws.Range("J57").Value = 17765.666
$ws.Range("L57").Value = 17765.666
$ws.Range("N57").Value = -19405.666
$ws.Range("H70").Value = 44806.848
$ws.Range("I70").Value = 46199.12
$ws.Range("K70").Value = 46199.12
$ws.Range("M70").Value = -45929.12
$ws.Range("H73").Value = 44806.848
$ws.Range("I73").Value = 46199.12
$ws.Range("K73").Value = 46199.12
$ws.Range("M73").Value = -45263.12
$ws.Range("H132").Value = 2431.8462
$ws.Range("I132").Value = 2117.35
$ws.Range("J132").Value = 3480.1667
$ws.Range("K132").Value = 6352.049999999999
$ws.Range("L132").Value = 10440.5001
$ws.Range("M132").Value = -3822.049999999999
$ws.Range("N132").Value = -15500.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4755.16
$ws.Range("I16").Value = 5285.591
$ws.Range("J16").Value = 865.3333
$ws.Range("K16").Value = 5285.591
$ws.Range("L16").Value = 865.3333
$ws.Range("M16").Value = -5115.591
$ws.Range("N16").Value = -1205.3333
$ws.Range("H61").Value = 2483.45
$ws.Range("I61").Value = 1999.4
$ws.Range("J61").Value = 2967.5
$ws.Range("K61").Value = 1999.4
$ws.Range("L61").Value = 2967.5
$ws.Range("M61").Value = -1797.4
$ws.Range("N61").Value = -3371.5
$ws.Range("H93").Value = 4000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 4000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -6496
$ws.Range("H113").Value = 2483.45
$ws.Range("I113").Value = 1999.4
$ws.Range("J113").Value = 2967.5
$ws.Range("K113").Value = 1999.4
$ws.Range("L113").Value = 2967.5
$ws.Range("M113").Value = 170.5999999999999
$ws.Range("N113").Value = -7307.5
$ws.Range("H122").Value = 3670.6667
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4097.4546
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 12292.3638
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -17192.3638
$ws.Range("H132").Value = 2710.0908
$ws.Range("I132").Value = 1601.1428
$ws.Range("J132").Value = 4650.75
$ws.Range("K132").Value = 4803.428400000001
$ws.Range("L132").Value = 13952.25
$ws.Range("M132").Value = -2273.428400000001
$ws.Range("N132").Value = -19012.25
$ws.Range("H133").Value = 40923.184
$ws.Range("J133").Value = 40923.184
$ws.Range("L133").Value = 40923.184
$ws.Range("N133").Value = -45983.184
